# Updates "想去人数" (F) / "最低票价" (G) counts across the four sheets
# of the 上海-漫展信息 workbook, per the gh-pages regeneration commit.

$wb = $excel.ActiveWorkbook

# sheetName -> { cellRef -> newValue }
$updates = @{
    "展览" = @{
        "F5"  = 203
        "F6"  = 1026
        "F9"  = 762
        "F10" = 1461
        "F11" = 1048
        "F12" = 784
        "G12" = 68
        "F13" = 792
        "F15" = 605
        "F17" = 696
        "F23" = 5350
        "F26" = 2477
        "F27" = 5936
        "F29" = 1030
        "F30" = 608
        "F33" = 1063
        "F35" = 76
        "F37" = 724
        "F42" = 14
        "F44" = 46
        "F47" = 572
    }
    "演出" = @{
        "F5"  = 2094
        "F6"  = 490
        "G10" = 118
        "F24" = 1713
        "F37" = 174
        "F38" = 910
        "F39" = 501
        "F41" = 38
        "F47" = 18
    }
    "本地生活" = @{
        "F5" = 780
        "F6" = 438
        "F7" = 250
    }
    "全部类型" = @{
        "F6"  = 438
        "F7"  = 250
        "F8"  = 250
        "F9"  = 203
        "F10" = 1026
        "F13" = 490
        "F14" = 762
        "F15" = 1461
        "F18" = 1048
        "F19" = 784
        "G19" = 68
        "G20" = 118
        "F21" = 792
        "F23" = 605
        "F24" = 696
        "F30" = 5350
        "F33" = 2477
        "F34" = 5936
        "F35" = 1030
        "F36" = 1713
        "F37" = 608
        "F39" = 1063
        "F40" = 76
        "F41" = 724
        "F45" = 910
        "F46" = 501
        "F47" = 46
        "F48" = 38
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
